$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    Write-Output ("$i : " + ('{0:X6}' -f $tcs.Item($i).RGB))
}
